$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 3020
$wsExhibit.Range("F5").Value = 82

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 3020
$wsAll.Range("F10").Value = 82
